$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
# A8 holds "Volume 29   Number  48" as rich-text runs; only the trailing
# volume-number run ("48" -> "49") changes.
$ws.Range("A8").Characters(21, 2).Text = "49"

# C9 holds "Report Covering the Week  11/28/2022  Through  12/4/2022" as
# rich-text runs; the two date runs change. The second run's start shifts
# left by one character once the (10-char -> 9-char) first date is replaced.
$ws.Range("C9").Characters(27, 10).Text = "12/5/2022"
$ws.Range("C9").Characters(47, 9).Text = "12/11/2022"

# --- Precinct crime-stat table refresh (rows 16-30) ---

# Row 16
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -53.846153846153
$ws.Range("I16").Value = 148
$ws.Range("J16").Value = 136
$ws.Range("K16").Value = 8.823529411764
$ws.Range("L16").Value = 11.278195488721
$ws.Range("M16").Value = -20
$ws.Range("N16").Value = -82.629107981220

# Row 17
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = -16.666666666666
$ws.Range("I17").Value = 177
$ws.Range("J17").Value = 176
$ws.Range("K17").Value = 0.568181818181
$ws.Range("L17").Value = 18.791946308724
$ws.Range("M17").Value = 77
$ws.Range("N17").Value = -60.753880266075

# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("I18").Value = 173
$ws.Range("J18").Value = 127
$ws.Range("K18").Value = 36.220472440944
$ws.Range("L18").Value = -9.424083769633
$ws.Range("M18").Value = 40.650406504065
$ws.Range("N18").Value = -83.039215686274

# Row 19
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 116.666666666667
$ws.Range("F19").Value = 43
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = 2.380952380952
$ws.Range("I19").Value = 504
$ws.Range("J19").Value = 451
$ws.Range("K19").Value = 11.751662971175
$ws.Range("L19").Value = 50.447761194029
$ws.Range("M19").Value = -1.5625
$ws.Range("N19").Value = -52.676056338028

# Row 20
$ws.Range("F20").Value = 4
$ws.Range("H20").Value = -50
$ws.Range("I20").Value = 83
$ws.Range("K20").Value = -1.190476190476
$ws.Range("L20").Value = 13.698630136986
$ws.Range("M20").Value = 102.439024390244
$ws.Range("N20").Value = -91.235480464625

# Row 21
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = 57.142857142857
$ws.Range("F21").Value = 83
$ws.Range("G21").Value = 95
$ws.Range("H21").Value = -12.631578947368
$ws.Range("I21").Value = 1095
$ws.Range("J21").Value = 989
$ws.Range("K21").Value = 10.717896865520
$ws.Range("L21").Value = 22.620380739081
$ws.Range("M21").Value = 12.422997946611
$ws.Range("N21").Value = -75.062628102937

# Row 22
$ws.Range("C22").Value = 2
$ws.Range("F22").Value = 4
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 32
$ws.Range("K22").Value = 23.076923076923
$ws.Range("L22").Value = 39.130434782608
$ws.Range("M22").Value = 10.344827586206

# Row 23 - C/D/E switch from numeric counts to the sheet's "no data" markers
# ("0" / "***.*"), reusing the same shared strings used elsewhere (rows 14,
# 20, 26, ...).
$ws.Range("C23").Value = "0"
$ws.Range("D23").Value = "0"
$ws.Range("E23").Value = "***.*"
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 14.285714285714
$ws.Range("I23").Value = 118
$ws.Range("K23").Value = 13.461538461538
$ws.Range("L23").Value = 28.260869565217
$ws.Range("M23").Value = 59.459459459459

# Row 24
$ws.Range("C24").Value = 36
$ws.Range("D24").Value = 35
$ws.Range("E24").Value = 2.857142857142
$ws.Range("F24").Value = 130
$ws.Range("G24").Value = 165
$ws.Range("H24").Value = -21.212121212121
$ws.Range("I24").Value = 1772
$ws.Range("J24").Value = 1220
$ws.Range("K24").Value = 45.245901639344
$ws.Range("L24").Value = 32.535527299925
$ws.Range("M24").Value = 74.409448818897

# Row 25
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -62.5
$ws.Range("F25").Value = 18
$ws.Range("G25").Value = 25
$ws.Range("H25").Value = -28
$ws.Range("I25").Value = 288
$ws.Range("J25").Value = 261
$ws.Range("K25").Value = 10.344827586206
$ws.Range("L25").Value = 20.502092050209
$ws.Range("M25").Value = -7.987220447284

# Row 26
$ws.Range("F26").Value = 1
$ws.Range("H26").Value = -50

# Row 27 - D/E switch from the "no data" markers to real numbers
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = -25
$ws.Range("J27").Value = 52
$ws.Range("K27").Value = -1.923076923076
$ws.Range("L27").Value = 30.769230769230

# Row 30
$ws.Range("L30").Value = 433.333333333333
